$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1
$ws.Range("H1").Value = "Save"

# Match the style of the existing headers (copy format from G1, the "sum" header)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add the value for the new column in H2
$ws.Range("H2").Value = 0
